$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column M (the "2019" data column). This shifts the following
# "Notes" column (N) left into M, matching the source edit exactly.
$ws.Columns("M").Delete()

# Restore the selection to the cell the author ended up on.
$ws.Range("A19").Select()
